# Update countries & provincias Spain
#
# This script updates the COVID-19 statistics table on sheet "Pais" to a
# newer data pull. Most updates are simple value refreshes, but several
# countries changed relative rank (because their "Casos totales" - column B -
# grew past a neighbour), so the row that a given country occupies in the
# sorted-by-B-descending table changes too. We therefore write both the
# country label (column A) and the data (columns B-H) for every row whose
# content changed, rather than only the numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos (value refresh only)
$ws.Range("B4").Value = 6838553
$ws.Range("C4").Value = 10252
$ws.Range("D4").Value = 4120574
$ws.Range("E4").Value = 2516348
$ws.Range("G4").Value = 283
$ws.Range("H4").Value = 201631

# Row 17: Reino Unido (value refresh only)
$ws.Range("B17").Value = 381614
$ws.Range("C17").Value = 3395
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 41705

# Row 25: Alemania (value refresh only)
$ws.Range("B25").Value = 267512
$ws.Range("C25").Value = 647
$ws.Range("E25").Value = 18961

# Row 29: Canada (value refresh only)
$ws.Range("B29").Value = 140040
$ws.Range("C29").Value = 293
$ws.Range("D29").Value = 122452
$ws.Range("E29").Value = 8392
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 9196

# Rows 51-52: Portugal overtakes Etiopia (rank swap)
$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 66396
$ws.Range("C51").Value = 770
$ws.Range("D51").Value = 44794
$ws.Range("E51").Value = 19714
$ws.Range("G51").Value = 10
$ws.Range("H51").Value = 1888

$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 66224
$ws.Range("D52").Value = 26665
$ws.Range("E52").Value = 38514
$ws.Range("H52").Value = 1045

# Row 57: Singapur (value refresh only)
$ws.Range("D57").Value = 57039
$ws.Range("E57").Value = 466

# Row 65: Moldavia (value refresh only)
$ws.Range("B65").Value = 44983
$ws.Range("C65").Value = 622
$ws.Range("D65").Value = 33239
$ws.Range("E65").Value = 10574
$ws.Range("G65").Value = 11
$ws.Range("H65").Value = 1170

# Rows 126-128: Jordania overtakes Eslovenia and Birmania (rank swap)
$ws.Range("A126").Value = "Jordania"
$ws.Range("B126").Value = 4131
$ws.Range("C126").Value = 279
$ws.Range("D126").Value = 2415
$ws.Range("E126").Value = 1690
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 26

$ws.Range("A127").Value = "Eslovenia"
$ws.Range("B127").Value = 4058
$ws.Range("C127").Value = 104
$ws.Range("D127").Value = 2897
$ws.Range("E127").Value = 1025
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 136

$ws.Range("A128").Value = "Birmania"
$ws.Range("B128").Value = 4043
$ws.Range("C128").Value = 222
$ws.Range("D128").Value = 944
$ws.Range("E128").Value = 3053
$ws.Range("G128").Value = 6
$ws.Range("H128").Value = 46

# Rows 214-215: Montserrat overtakes Islas Malvinas (tie-break rank swap)
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Update the "last updated" timestamp shown at the top of the sheet
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 17:32"
